# Set column D ("¿Resuelto? (✓/✗)") value to "si" for exercises 43 to 48
# (Excel rows 44 through 49), leaving the rest of the sheet untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D44:D49").Value = "si"

# Restore the view state (scroll position / selection) to match what the
# author left the sheet on after editing. Best-effort: some hosts may not
# expose window scrolling, so don't let that abort the data edit above.
try {
    $ws.Activate()
    $ws.Range("C49").Select()
    $excel.ActiveWindow.ScrollRow = 34
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # View-state restore is cosmetic; ignore failures here.
}
